$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append new user rows to the "Munka1" (Sheet1) user table.
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "tibi"
$ws.Range("C5").Value = "tibi"

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "ferixx"
$ws.Range("C6").Value = "xx"

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "asd"
$ws.Range("C7").Value = "asd"

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "newuser"
$ws.Range("C8").Value = "iiii"
